# Insert a new data row for Femacal de La Calera / Zapallo italiano just
# above the existing row 471, shifting rows 471-494 down to 472-495, and
# populate the newly inserted row with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 471..494 down by one row (creates a new, mostly blank row 471)
$ws.Rows(471).Insert()

# Populate the new row 471 with the new record
$ws.Cells.Item(471, 1).Value = 3
$ws.Cells.Item(471, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(471, 3).Value = "Coquimbo"
$ws.Cells.Item(471, 4).Value = 44753
$ws.Cells.Item(471, 5).Value = 5
$ws.Cells.Item(471, 6).Value = 100112032
$ws.Cells.Item(471, 7).Value = "Zapallo italiano"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 190
$ws.Cells.Item(471, 11).Value = 9000
$ws.Cells.Item(471, 12).Value = 10000
$ws.Cells.Item(471, 13).Value = 9487
$ws.Cells.Item(471, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(471, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(471, 16).Value = 136
$ws.Cells.Item(471, 17).Value = 70
$ws.Cells.Item(471, 18).Value = "Hortaliza"
